$d = $word.ActiveDocument

# 1. Title: "DOVE TROVARE SMARTCASH?" -> "DOVE OTTENERE SMARTCASH?"
$d.Content.Find.Execute("DOVE TROVARE SMARTCASH?", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DOVE OTTENERE SMARTCASH?", 2)

# 2. Remove "quindi " before "essere un pochino"
$d.Content.Find.Execute("procurarsene un pò potrebbe quindi essere un pochino", $true, $false, $false, $false, $false,
                         $true, 1, $false, "procurarsene un pò potrebbe essere un pochino", 2)

# 3. "PORTAFOGLIO" -> "WALLET"
$d.Content.Find.Execute("PORTAFOGLIO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "WALLET", 2)

# 4. "Per il supporto sul mining cortesemente raggiungi " -> "Per il supporto sul mining accedi su "
$d.Content.Find.Execute("Per il supporto sul mining cortesemente raggiungi ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Per il supporto sul mining accedi su ", 2)
